$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data: new values for columns C,D,E,F,G,H,I,J,N across rows 2-25
$data = @{}
$data[2] = @(4.067937038739815, 7.775327537238806, 12.79089778093609, 43.55292226603626, 55.81828364157803, 19.41613829829461, 30.49359561336956, 10.09421402311757, 18.99769123883691)
$data[3] = @(4.079197275270559, 7.783950329519477, 12.79621159339102, 42.80830285766464, 54.23993323042032, 19.23223933023548, 29.98868480780467, 10.0996582748656, 18.40031508502705)
$data[4] = @(4.086823181783126, 7.78959301488488, 12.80222834265451, 42.36012949402318, 53.26890903784871, 19.12460663806837, 29.68522153962684, 10.10513802733391, 18.02485520896364)
$data[5] = @(4.090109392622828, 7.791980221878919, 12.80537095394674, 42.18000477163218, 52.87339574061173, 19.08211190537825, 29.56336527187998, 10.10790694008695, 17.86990355188116)
$data[6] = @(4.090665838468388, 7.792381922476633, 12.80593444346415, 42.15025368690427, 52.80775238547449, 19.07513927552335, 29.54324488440697, 10.10839903397151, 17.84406337566596)
$data[7] = @(4.086866778263603, 7.789624853964821, 12.80226793078284, 42.35768979431207, 53.26357334810292, 19.1240279581947, 29.68357061698875, 10.10517320223392, 18.02277304766463)
$data[8] = @(4.071671495189759, 7.778228512134253, 12.79215737886158, 43.29444933870976, 55.2748563996356, 19.35165752383271, 30.3182387230674, 10.09564683369204, 18.79364780656866)
$data[9] = @(4.047544973677527, 7.758634923494901, 12.79425332100533, 45.19235793348977, 59.17608532848412, 19.83812848513048, 31.60769287801812, 10.09398137337363, 20.2273683202997)
$data[10] = @(4.033306825936784, 7.745907050268788, 12.80924056609112, 46.60957261525966, 61.98165221265616, 20.21737083869742, 32.57281943039996, 10.10320331252063, 21.22223697909767)
$data[11] = @(4.027593510409077, 7.740476480735686, 12.81899056674729, 47.25628024516141, 63.23832606388049, 20.39404390176285, 33.0137342405848, 10.10967815739478, 21.66018057919901)
$data[12] = @(4.025540408111198, 7.738471567625872, 12.8231046955808, 47.50122072253081, 63.71089757901808, 20.46149083548896, 33.18080441632883, 10.11245835723047, 21.82377585682186)
$data[13] = @(4.025977660591148, 7.738901072355381, 12.82219987252583, 47.4484701162484, 63.60927505002351, 20.44694155215288, 33.14482073616266, 10.11184498395798, 21.78864458690801)
$data[14] = @(4.027422384339199, 7.740310503551906, 12.81932057978252, 47.27643207546983, 63.2772736443035, 20.39958218944841, 33.02747801478522, 10.10990030410748, 21.67368539489659)
$data[15] = @(4.028321715129525, 7.74118052647524, 12.81761189429875, 47.171052790493, 63.07346929284495, 20.37064257097967, 32.95561114459639, 10.10875189896328, 21.60297336126124)
$data[16] = @(4.033695622037798, 7.746269169429544, 12.80866242263741, 46.56732825125147, 61.89908879863994, 20.20590394254982, 32.54402816467353, 10.10282608747741, 21.19330956972086)
$data[17] = @(4.037188354973475, 7.749482826417784, 12.8039235813655, 46.19730219973583, 61.17326632242001, 20.10586883792672, 32.29189684943368, 10.09977510984827, 20.93814219015166)
$data[18] = @(4.039269128899833, 7.751365074753176, 12.80147391798385, 45.98467618138552, 60.75397454483106, 20.04872548346213, 32.14706349118212, 10.09823483244172, 20.79000725568362)
$data[19] = @(4.039985965060463, 7.752008188215047, 12.80069188997182, 45.91272709304162, 60.61171286734237, 20.02944701854527, 32.09806242153538, 10.09775015397219, 20.73962067985785)
$data[20] = @(4.036809107536405, 7.749137226204421, 12.80439947130398, 46.23667286139244, 61.25072317052038, 20.11647732707275, 32.31871862397813, 10.10007768078277, 20.96544799483449)
$data[21] = @(4.026995032819365, 7.73989512227204, 12.82015484309059, 47.32696445070297, 63.37488382026352, 20.41347842282217, 33.06194283953486, 10.11046259045562, 21.70751365554064)
$data[22] = @(4.021224779697764, 7.734155135375105, 12.83291167125705, 48.03969709016903, 64.74369508337041, 20.6107357861555, 33.54822457537121, 10.11916348416585, 22.17935961385674)
$data[23] = @(4.024245366627381, 7.737191250719981, 12.82587799992599, 47.65936106480507, 64.01506375992464, 20.50518514818209, 33.28869005859318, 10.11434443829059, 21.92877110911181)
$data[24] = @(4.036980338731152, 7.749293364046929, 12.80418346567895, 46.21887303382024, 61.21571114275402, 20.11168007841129, 32.30659211250465, 10.09994022265642, 20.95310750188673)
$data[25] = @(4.053461589684665, 7.76364183937358, 12.79132965375935, 44.67388760418577, 58.12901065079915, 19.70249413479848, 31.255051104559, 10.09260199049056, 19.84905939529497)

$cols = @(3, 4, 5, 6, 7, 8, 9, 10, 14)  # C, D, E, F, G, H, I, J, N

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $ws.Cells.Item($row, $cols[$i]).Value = [double]$vals[$i]
    }
}